$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55
$ws.Cells.Item(55, 2).Value = 6830657
$ws.Cells.Item(55, 5).Value = "Chivas Guadalajara Women"
$ws.Cells.Item(55, 6).Value = "Monterrey Women"
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 1
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = "H"
$ws.Cells.Item(55, 12).Value = 2.375
$ws.Cells.Item(55, 13).Value = 3.4
$ws.Cells.Item(55, 14).Value = 2.5
$ws.Cells.Item(55, 15).Value = 2.15
$ws.Cells.Item(55, 16).Value = 3.5
$ws.Cells.Item(55, 17).Value = 3
$ws.Cells.Item(55, 18).Value = -0.25
$ws.Cells.Item(55, 19).Value = 1.875
$ws.Cells.Item(55, 20).Value = 1.925
$ws.Cells.Item(55, 21).Value = 2.75
$ws.Cells.Item(55, 22).Value = 1.75
$ws.Cells.Item(55, 23).Value = 1.95
$ws.Cells.Item(55, 24).Value = 1.15
$ws.Cells.Item(55, 25).Value = -1
$ws.Cells.Item(55, 26).Value = -1
$ws.Cells.Item(55, 27).Value = 0.875
$ws.Cells.Item(55, 28).Value = -1
$ws.Cells.Item(55, 29).Value = -1
$ws.Cells.Item(55, 30).Value = 0.95
# Row 56
$ws.Cells.Item(56, 2).Value = 6830656
$ws.Cells.Item(56, 5).Value = "Mazatlan FC Women"
$ws.Cells.Item(56, 6).Value = "Club Necaxa Women"
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 2
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 1
$ws.Cells.Item(56, 11).Value = "A"
$ws.Cells.Item(56, 12).Value = 3.5
$ws.Cells.Item(56, 13).Value = 3.5
$ws.Cells.Item(56, 14).Value = 1.833
$ws.Cells.Item(56, 15).Value = 4.333
$ws.Cells.Item(56, 16).Value = 3.75
$ws.Cells.Item(56, 17).Value = 1.727
$ws.Cells.Item(56, 18).Value = 0.75
$ws.Cells.Item(56, 19).Value = 1.85
$ws.Cells.Item(56, 20).Value = 1.95
$ws.Cells.Item(56, 21).Value = 2.75
$ws.Cells.Item(56, 22).Value = 1.975
$ws.Cells.Item(56, 23).Value = 1.825
$ws.Cells.Item(56, 24).Value = -1
$ws.Cells.Item(56, 25).Value = -1
$ws.Cells.Item(56, 26).Value = 0.7270000000000001
$ws.Cells.Item(56, 27).Value = -1
$ws.Cells.Item(56, 28).Value = 0.95
$ws.Cells.Item(56, 29).Value = -1
$ws.Cells.Item(56, 30).Value = 0.825
# Row 229
$ws.Cells.Item(229, 2).Value = 7645712
$ws.Cells.Item(229, 5).Value = "Chivas Guadalajara Women"
$ws.Cells.Item(229, 6).Value = "Santos Laguna Women"
$ws.Cells.Item(229, 7).Value = 10
$ws.Cells.Item(229, 8).Value = 2
$ws.Cells.Item(229, 9).Value = 6
$ws.Cells.Item(229, 10).Value = 0
$ws.Cells.Item(229, 11).Value = "H"
$ws.Cells.Item(229, 12).Value = 1.1
$ws.Cells.Item(229, 13).Value = 8
$ws.Cells.Item(229, 14).Value = 13
$ws.Cells.Item(229, 15).Value = 1.03
$ws.Cells.Item(229, 16).Value = 17
$ws.Cells.Item(229, 17).Value = 41
$ws.Cells.Item(229, 18).Value = -3.75
$ws.Cells.Item(229, 19).Value = 1.775
$ws.Cells.Item(229, 20).Value = 1.925
$ws.Cells.Item(229, 21).Value = 4.75
$ws.Cells.Item(229, 22).Value = 1.9
$ws.Cells.Item(229, 23).Value = 1.9
$ws.Cells.Item(229, 24).Value = 0.03000000000000003
$ws.Cells.Item(229, 25).Value = -1
$ws.Cells.Item(229, 26).Value = -1
$ws.Cells.Item(229, 27).Value = 0.7749999999999999
$ws.Cells.Item(229, 28).Value = -1
$ws.Cells.Item(229, 29).Value = 0.8999999999999999
$ws.Cells.Item(229, 30).Value = -1
# Row 230
$ws.Cells.Item(230, 2).Value = 7645781
$ws.Cells.Item(230, 5).Value = "Atletico San Luis Women"
$ws.Cells.Item(230, 6).Value = "Cruz Azul Women"
$ws.Cells.Item(230, 7).Value = 3
$ws.Cells.Item(230, 8).Value = 5
$ws.Cells.Item(230, 9).Value = 1
$ws.Cells.Item(230, 10).Value = 3
$ws.Cells.Item(230, 11).Value = "A"
$ws.Cells.Item(230, 12).Value = 2.1
$ws.Cells.Item(230, 13).Value = 3.6
$ws.Cells.Item(230, 14).Value = 2.8
$ws.Cells.Item(230, 15).Value = 2.45
$ws.Cells.Item(230, 16).Value = 3.6
$ws.Cells.Item(230, 17).Value = 2.375
$ws.Cells.Item(230, 18).Value = 0
$ws.Cells.Item(230, 19).Value = 1.95
$ws.Cells.Item(230, 20).Value = 1.85
$ws.Cells.Item(230, 21).Value = 2.75
$ws.Cells.Item(230, 22).Value = 1.75
$ws.Cells.Item(230, 23).Value = 1.95
$ws.Cells.Item(230, 24).Value = -1
$ws.Cells.Item(230, 25).Value = -1
$ws.Cells.Item(230, 26).Value = 1.375
$ws.Cells.Item(230, 27).Value = -1
$ws.Cells.Item(230, 28).Value = 0.8500000000000001
$ws.Cells.Item(230, 29).Value = 0.75
$ws.Cells.Item(230, 30).Value = -1
# Row 231
$ws.Cells.Item(231, 2).Value = 7645713
$ws.Cells.Item(231, 5).Value = "Toluca Women"
$ws.Cells.Item(231, 6).Value = "Tigres UANL Women"
$ws.Cells.Item(231, 7).Value = 0
$ws.Cells.Item(231, 8).Value = 7
$ws.Cells.Item(231, 9).Value = 0
$ws.Cells.Item(231, 10).Value = 3
$ws.Cells.Item(231, 11).Value = "A"
$ws.Cells.Item(231, 12).Value = 9
$ws.Cells.Item(231, 13).Value = 7
$ws.Cells.Item(231, 14).Value = 1.166
$ws.Cells.Item(231, 15).Value = 5.25
$ws.Cells.Item(231, 16).Value = 4.75
$ws.Cells.Item(231, 17).Value = 1.4
$ws.Cells.Item(231, 18).Value = 1.25
$ws.Cells.Item(231, 19).Value = 1.925
$ws.Cells.Item(231, 20).Value = 1.875
$ws.Cells.Item(231, 21).Value = 3
$ws.Cells.Item(231, 22).Value = 1.75
$ws.Cells.Item(231, 23).Value = 1.95
$ws.Cells.Item(231, 24).Value = -1
$ws.Cells.Item(231, 25).Value = -1
$ws.Cells.Item(231, 26).Value = 0.3999999999999999
$ws.Cells.Item(231, 27).Value = -1
$ws.Cells.Item(231, 28).Value = 0.875
$ws.Cells.Item(231, 29).Value = 0.75
$ws.Cells.Item(231, 30).Value = -1
# Row 245
$ws.Cells.Item(245, 2).Value = 7645794
$ws.Cells.Item(245, 5).Value = "Toluca Women"
$ws.Cells.Item(245, 6).Value = "Club America Women"
$ws.Cells.Item(245, 7).Value = 3
$ws.Cells.Item(245, 8).Value = 0
$ws.Cells.Item(245, 9).Value = 1
$ws.Cells.Item(245, 10).Value = 0
$ws.Cells.Item(245, 11).Value = "H"
$ws.Cells.Item(245, 12).Value = 6
$ws.Cells.Item(245, 13).Value = 5
$ws.Cells.Item(245, 14).Value = 1.333
$ws.Cells.Item(245, 15).Value = 5.75
$ws.Cells.Item(245, 16).Value = 5
$ws.Cells.Item(245, 17).Value = 1.363
$ws.Cells.Item(245, 18).Value = 1.5
$ws.Cells.Item(245, 19).Value = 1.825
$ws.Cells.Item(245, 20).Value = 1.975
$ws.Cells.Item(245, 21).Value = 3.25
$ws.Cells.Item(245, 22).Value = 1.825
$ws.Cells.Item(245, 23).Value = 1.975
$ws.Cells.Item(245, 24).Value = 4.75
$ws.Cells.Item(245, 25).Value = -1
$ws.Cells.Item(245, 26).Value = -1
$ws.Cells.Item(245, 27).Value = 0.825
$ws.Cells.Item(245, 28).Value = -1
$ws.Cells.Item(245, 29).Value = -0.5
$ws.Cells.Item(245, 30).Value = 0.4875
# Row 246
$ws.Cells.Item(246, 2).Value = 7645793
$ws.Cells.Item(246, 5).Value = "Queretaro Women"
$ws.Cells.Item(246, 6).Value = "Cruz Azul Women"
$ws.Cells.Item(246, 7).Value = 3
$ws.Cells.Item(246, 8).Value = 0
$ws.Cells.Item(246, 9).Value = 1
$ws.Cells.Item(246, 10).Value = 0
$ws.Cells.Item(246, 11).Value = "H"
$ws.Cells.Item(246, 12).Value = 1.909
$ws.Cells.Item(246, 13).Value = 3.5
$ws.Cells.Item(246, 14).Value = 3.3
$ws.Cells.Item(246, 15).Value = 2.1
$ws.Cells.Item(246, 16).Value = 3.4
$ws.Cells.Item(246, 17).Value = 2.875
$ws.Cells.Item(246, 18).Value = -0.25
$ws.Cells.Item(246, 19).Value = 1.875
$ws.Cells.Item(246, 20).Value = 1.925
$ws.Cells.Item(246, 21).Value = 2.75
$ws.Cells.Item(246, 22).Value = 2
$ws.Cells.Item(246, 23).Value = 1.8
$ws.Cells.Item(246, 24).Value = 1.1
$ws.Cells.Item(246, 25).Value = -1
$ws.Cells.Item(246, 26).Value = -1
$ws.Cells.Item(246, 27).Value = 0.875
$ws.Cells.Item(246, 28).Value = -1
$ws.Cells.Item(246, 29).Value = 0.5
$ws.Cells.Item(246, 30).Value = -0.5
# Row 248
$ws.Cells.Item(248, 2).Value = 7645719
$ws.Cells.Item(248, 5).Value = "Monterrey Women"
$ws.Cells.Item(248, 6).Value = "Santos Laguna Women"
$ws.Cells.Item(248, 7).Value = 6
$ws.Cells.Item(248, 8).Value = 0
$ws.Cells.Item(248, 9).Value = 4
$ws.Cells.Item(248, 10).Value = 0
$ws.Cells.Item(248, 11).Value = "H"
$ws.Cells.Item(248, 12).Value = 1.025
$ws.Cells.Item(248, 13).Value = 15
$ws.Cells.Item(248, 14).Value = 34
$ws.Cells.Item(248, 15).Value = 1.01
$ws.Cells.Item(248, 16).Value = 34
$ws.Cells.Item(248, 17).Value = 67
$ws.Cells.Item(248, 18).Value = -4.75
$ws.Cells.Item(248, 19).Value = 1.775
$ws.Cells.Item(248, 20).Value = 1.925
$ws.Cells.Item(248, 21).Value = 5.75
$ws.Cells.Item(248, 22).Value = 1.85
$ws.Cells.Item(248, 23).Value = 1.95
$ws.Cells.Item(248, 24).Value = 0.01000000000000001
$ws.Cells.Item(248, 25).Value = -1
$ws.Cells.Item(248, 26).Value = -1
$ws.Cells.Item(248, 27).Value = 0.7749999999999999
$ws.Cells.Item(248, 28).Value = -1
$ws.Cells.Item(248, 29).Value = 0.425
$ws.Cells.Item(248, 30).Value = -0.5
# Row 249
$ws.Cells.Item(249, 2).Value = 7645796
$ws.Cells.Item(249, 5).Value = "Chivas Guadalajara Women"
$ws.Cells.Item(249, 6).Value = "Club Necaxa Women"
$ws.Cells.Item(249, 7).Value = 4
$ws.Cells.Item(249, 8).Value = 0
$ws.Cells.Item(249, 9).Value = 3
$ws.Cells.Item(249, 10).Value = 0
$ws.Cells.Item(249, 11).Value = "H"
$ws.Cells.Item(249, 12).Value = 1.055
$ws.Cells.Item(249, 13).Value = 10
$ws.Cells.Item(249, 14).Value = 21
$ws.Cells.Item(249, 15).Value = 1.062
$ws.Cells.Item(249, 16).Value = 11
$ws.Cells.Item(249, 17).Value = 29
$ws.Cells.Item(249, 18).Value = -3
$ws.Cells.Item(249, 19).Value = 1.85
$ws.Cells.Item(249, 20).Value = 1.95
$ws.Cells.Item(249, 21).Value = 4
$ws.Cells.Item(249, 22).Value = 1.8
$ws.Cells.Item(249, 23).Value = 2
$ws.Cells.Item(249, 24).Value = 0.06200000000000006
$ws.Cells.Item(249, 25).Value = -1
$ws.Cells.Item(249, 26).Value = -1
$ws.Cells.Item(249, 27).Value = 0.8500000000000001
$ws.Cells.Item(249, 28).Value = -1
$ws.Cells.Item(249, 29).Value = 0
$ws.Cells.Item(249, 30).Value = 0
# Row 271
$ws.Cells.Item(271, 2).Value = 7645809
$ws.Cells.Item(271, 5).Value = "Mazatlan FC Women"
$ws.Cells.Item(271, 6).Value = "Queretaro Women"
$ws.Cells.Item(271, 7).Value = 2
$ws.Cells.Item(271, 8).Value = 2
$ws.Cells.Item(271, 9).Value = 0
$ws.Cells.Item(271, 10).Value = 1
$ws.Cells.Item(271, 11).Value = "D"
$ws.Cells.Item(271, 12).Value = 4.75
$ws.Cells.Item(271, 13).Value = 4
$ws.Cells.Item(271, 14).Value = 1.533
$ws.Cells.Item(271, 15).Value = 4.75
$ws.Cells.Item(271, 16).Value = 3.75
$ws.Cells.Item(271, 17).Value = 1.571
$ws.Cells.Item(271, 18).Value = 1
$ws.Cells.Item(271, 19).Value = 1.775
$ws.Cells.Item(271, 20).Value = 2.025
$ws.Cells.Item(271, 21).Value = 2.75
$ws.Cells.Item(271, 22).Value = 1.95
$ws.Cells.Item(271, 23).Value = 1.85
$ws.Cells.Item(271, 24).Value = -1
$ws.Cells.Item(271, 25).Value = 2.75
$ws.Cells.Item(271, 26).Value = -1
$ws.Cells.Item(271, 27).Value = 0.7749999999999999
$ws.Cells.Item(271, 28).Value = -1
$ws.Cells.Item(271, 29).Value = 0.95
$ws.Cells.Item(271, 30).Value = -1
# Row 272
$ws.Cells.Item(272, 2).Value = 7645812
$ws.Cells.Item(272, 5).Value = "Atlas Women"
$ws.Cells.Item(272, 6).Value = "Monterrey Women"
$ws.Cells.Item(272, 7).Value = 0
$ws.Cells.Item(272, 8).Value = 1
$ws.Cells.Item(272, 9).Value = 0
$ws.Cells.Item(272, 10).Value = 1
$ws.Cells.Item(272, 11).Value = "A"
$ws.Cells.Item(272, 12).Value = 6.5
$ws.Cells.Item(272, 13).Value = 5
$ws.Cells.Item(272, 14).Value = 1.333
$ws.Cells.Item(272, 15).Value = 6.5
$ws.Cells.Item(272, 16).Value = 5
$ws.Cells.Item(272, 17).Value = 1.3
$ws.Cells.Item(272, 18).Value = 1.5
$ws.Cells.Item(272, 19).Value = 1.875
$ws.Cells.Item(272, 20).Value = 1.925
$ws.Cells.Item(272, 21).Value = 3.25
$ws.Cells.Item(272, 22).Value = 2
$ws.Cells.Item(272, 23).Value = 1.8
$ws.Cells.Item(272, 24).Value = -1
$ws.Cells.Item(272, 25).Value = -1
$ws.Cells.Item(272, 26).Value = 0.3
$ws.Cells.Item(272, 27).Value = 0.875
$ws.Cells.Item(272, 28).Value = -1
$ws.Cells.Item(272, 29).Value = -1
$ws.Cells.Item(272, 30).Value = 0.8
# Row 310
$ws.Cells.Item(310, 2).Value = 7645737
$ws.Cells.Item(310, 5).Value = "Toluca Women"
$ws.Cells.Item(310, 6).Value = "Tijuana Women"
$ws.Cells.Item(310, 7).Value = 2
$ws.Cells.Item(310, 8).Value = 1
$ws.Cells.Item(310, 9).Value = 0
$ws.Cells.Item(310, 10).Value = 0
$ws.Cells.Item(310, 11).Value = "H"
$ws.Cells.Item(310, 12).Value = 2.3
$ws.Cells.Item(310, 13).Value = 3.75
$ws.Cells.Item(310, 14).Value = 2.4
$ws.Cells.Item(310, 15).Value = 2.2
$ws.Cells.Item(310, 16).Value = 3.8
$ws.Cells.Item(310, 17).Value = 2.625
$ws.Cells.Item(310, 18).Value = -0.25
$ws.Cells.Item(310, 19).Value = 2
$ws.Cells.Item(310, 20).Value = 1.8
$ws.Cells.Item(310, 21).Value = 3
$ws.Cells.Item(310, 22).Value = 1.8
$ws.Cells.Item(310, 23).Value = 2
$ws.Cells.Item(310, 24).Value = 1.2
$ws.Cells.Item(310, 25).Value = -1
$ws.Cells.Item(310, 26).Value = -1
$ws.Cells.Item(310, 27).Value = 1
$ws.Cells.Item(310, 28).Value = -1
$ws.Cells.Item(310, 29).Value = 0
$ws.Cells.Item(310, 30).Value = 0
# Row 311
$ws.Cells.Item(311, 2).Value = 7645837
$ws.Cells.Item(311, 5).Value = "Queretaro Women"
$ws.Cells.Item(311, 6).Value = "Atlas Women"
$ws.Cells.Item(311, 7).Value = 0
$ws.Cells.Item(311, 8).Value = 2
$ws.Cells.Item(311, 9).Value = 0
$ws.Cells.Item(311, 10).Value = 1
$ws.Cells.Item(311, 11).Value = "A"
$ws.Cells.Item(311, 12).Value = 1.8
$ws.Cells.Item(311, 13).Value = 3.75
$ws.Cells.Item(311, 14).Value = 3.5
$ws.Cells.Item(311, 15).Value = 2.05
$ws.Cells.Item(311, 16).Value = 3.4
$ws.Cells.Item(311, 17).Value = 3
$ws.Cells.Item(311, 18).Value = -0.25
$ws.Cells.Item(311, 19).Value = 1.825
$ws.Cells.Item(311, 20).Value = 1.975
$ws.Cells.Item(311, 21).Value = 2.5
$ws.Cells.Item(311, 22).Value = 1.8
$ws.Cells.Item(311, 23).Value = 2
$ws.Cells.Item(311, 24).Value = -1
$ws.Cells.Item(311, 25).Value = -1
$ws.Cells.Item(311, 26).Value = 2
$ws.Cells.Item(311, 27).Value = -1
$ws.Cells.Item(311, 28).Value = 0.9750000000000001
$ws.Cells.Item(311, 29).Value = -1
$ws.Cells.Item(311, 30).Value = 1
